$d = $word.ActiveDocument

$replacements = @(
    @("32×15=", "74×77="),
    @("48×96=", "22×25="),
    @("56×78=", "82×74="),
    @("60×35=", "44×64="),
    @("87×81=", "54×77="),
    @("91×98=", "26×54="),
    @("86×89=", "14×67="),
    @("99×88=", "50×49="),
    @("83×58=", "30×17="),
    @("49×80=", "82×44="),
    @("25×57=", "66×46="),
    @("50×66=", "26×43="),
    @("95×72=", "21×70="),
    @("76×92=", "84×42="),
    @("90×98=", "17×56="),
    @("31×84=", "93×90="),
    @("38×63=", "39×64="),
    @("49×96=", "89×85="),
    @("48×57=", "64×77="),
    @("24×27=", "53×98="),
    @("34×38=", "54×98="),
    @("33×21=", "90×27="),
    @("92×95=", "12×80="),
    @("23×42=", "89×97="),
    @("64×15=", "16×90=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
